# Fixed typo on slide
# Slide 2 ("Sorting an Array" / "Selection sort") is retitled:
#   Title:    "Sorting an Array" -> "Sortedness of an Array"
#   Subtitle: "Selection sort"   -> "Two definitions of "sorted""  (curly quotes)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Title shape: "Sorting an Array" -> "Sortedness of an Array" ---
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange

# Original text is "Sorting an Array":
#   chars 1-7   => "Sorting"
#   char  8     => " "
#   chars 9-16  => "an Array"

# 1) "Sorting" -> "Sortedness"
$titleRange.Characters(1, 7).Text = "Sortedness"

# 2) the single space that used to separate "Sorting" and "an" becomes " of "
#    (text is now "Sortedness an Array"; the separating space sits right after
#    the 10 characters of "Sortedness")
$titleRange.Characters(11, 1).Text = " of "

# --- Subtitle textbox: "Selection sort" -> 'Two definitions of "sorted"' ---
$leftQuote = [char]0x201C
$rightQuote = [char]0x201D
$subtitle = $s.Shapes.Item(3)
$subtitle.TextFrame.TextRange.Text = "Two definitions of " + $leftQuote + "sorted" + $rightQuote
